$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 282
$ws.Range("B282").Value = 6989704
$ws.Range("C282").Value = 'Serbia Prva Liga'
$ws.Range("D282").Value = 'Serbia Prva Liga'
$ws.Range("E282").Value = 45270.375
$ws.Range("F282").Value = 'Radnicki Sremska Mitrovica'
$ws.Range("G282").Value = 'Metalac Gornji'
$ws.Range("H282").Value = 0
$ws.Range("I282").Value = 1
$ws.Range("J282").Value = 'A'
$ws.Range("K282").Value = 1.666
$ws.Range("L282").Value = 3.4
$ws.Range("M282").Value = 4.5
$ws.Range("N282").Value = 2.05
$ws.Range("O282").Value = 3.2
$ws.Range("P282").Value = 3.3
$ws.Range("Q282").Value = -0.25
$ws.Range("R282").Value = 1.825
$ws.Range("S282").Value = 1.975
$ws.Range("T282").Value = 1.75
$ws.Range("U282").Value = 1.8
$ws.Range("V282").Value = 2
$ws.Range("W282").Value = -1
$ws.Range("X282").Value = -1
$ws.Range("Y282").Value = 2.3
$ws.Range("Z282").Value = -1
$ws.Range("AA282").Value = 0.9750000000000001
$ws.Range("AB282").Value = -1
$ws.Range("AC282").Value = 1

# Row 283
$ws.Range("B283").Value = 6989332
$ws.Range("C283").Value = 'Serbia Prva Liga'
$ws.Range("D283").Value = 'Serbia Prva Liga'
$ws.Range("E283").Value = 45270.375
$ws.Range("F283").Value = 'OFK Belgrade'
$ws.Range("G283").Value = 'FK Dubocica'
$ws.Range("H283").Value = 1
$ws.Range("I283").Value = 0
$ws.Range("J283").Value = 'H'
$ws.Range("K283").Value = 1.4
$ws.Range("L283").Value = 4
$ws.Range("M283").Value = 7
$ws.Range("N283").Value = 1.285
$ws.Range("O283").Value = 4.333
$ws.Range("P283").Value = 11
$ws.Range("Q283").Value = -1.5
$ws.Range("R283").Value = 1.85
$ws.Range("S283").Value = 1.95
$ws.Range("T283").Value = 2.5
$ws.Range("U283").Value = 1.95
$ws.Range("V283").Value = 1.85
$ws.Range("W283").Value = 0.2849999999999999
$ws.Range("X283").Value = -1
$ws.Range("Y283").Value = -1
$ws.Range("Z283").Value = -1
$ws.Range("AA283").Value = 0.95
$ws.Range("AB283").Value = -1
$ws.Range("AC283").Value = 0.8500000000000001

# Row 284
$ws.Range("B284").Value = 6989633
$ws.Range("C284").Value = 'Serbia Prva Liga'
$ws.Range("D284").Value = 'Serbia Prva Liga'
$ws.Range("E284").Value = 45270.375
$ws.Range("F284").Value = 'FK Indija'
$ws.Range("G284").Value = 'OFK Vrsac'
$ws.Range("H284").Value = 2
$ws.Range("I284").Value = 1
$ws.Range("J284").Value = 'H'
$ws.Range("K284").Value = 1.8
$ws.Range("L284").Value = 3.1
$ws.Range("M284").Value = 4.2
$ws.Range("N284").Value = 1.75
$ws.Range("O284").Value = 3.2
$ws.Range("P284").Value = 4.75
$ws.Range("Q284").Value = -0.75
$ws.Range("R284").Value = 1.95
$ws.Range("S284").Value = 1.75
$ws.Range("T284").Value = 2
$ws.Range("U284").Value = 1.975
$ws.Range("V284").Value = 1.725
$ws.Range("W284").Value = 0.75
$ws.Range("X284").Value = -1
$ws.Range("Y284").Value = -1
$ws.Range("Z284").Value = 0.475
$ws.Range("AA284").Value = -0.5
$ws.Range("AB284").Value = 0.9750000000000001
$ws.Range("AC284").Value = -1

# Row 285
$ws.Range("B285").Value = 6989530
$ws.Range("C285").Value = 'Serbia Prva Liga'
$ws.Range("D285").Value = 'Serbia Prva Liga'
$ws.Range("E285").Value = 45271.375
$ws.Range("F285").Value = 'FK Tekstilac Odzaci'
$ws.Range("G285").Value = 'FK Radnicki Beograd'
$ws.Range("H285").Value = 5
$ws.Range("I285").Value = 0
$ws.Range("J285").Value = 'H'
$ws.Range("K285").Value = 1.666
$ws.Range("L285").Value = 3.4
$ws.Range("M285").Value = 4.5
$ws.Range("N285").Value = 1.444
$ws.Range("O285").Value = 4.333
$ws.Range("P285").Value = 5.25
$ws.Range("Q285").Value = -1
$ws.Range("R285").Value = 1.875
$ws.Range("S285").Value = 1.925
$ws.Range("T285").Value = 2.5
$ws.Range("U285").Value = 1.85
$ws.Range("V285").Value = 1.95
$ws.Range("W285").Value = 0.444
$ws.Range("X285").Value = -1
$ws.Range("Y285").Value = -1
$ws.Range("Z285").Value = 0.875
$ws.Range("AA285").Value = -1
$ws.Range("AB285").Value = 0.8500000000000001
$ws.Range("AC285").Value = -1

# Row 286
$ws.Range("B286").Value = 6989531
$ws.Range("C286").Value = 'Serbia Prva Liga'
$ws.Range("D286").Value = 'Serbia Prva Liga'
$ws.Range("E286").Value = 45271.375
$ws.Range("F286").Value = 'FK Graficar Beograd'
$ws.Range("G286").Value = 'RFK Novi Sad 1921'
$ws.Range("H286").Value = 2
$ws.Range("I286").Value = 1
$ws.Range("J286").Value = 'H'
$ws.Range("K286").Value = 1.571
$ws.Range("L286").Value = 3.6
$ws.Range("M286").Value = 5
$ws.Range("N286").Value = 1.444
$ws.Range("O286").Value = 4.2
$ws.Range("P286").Value = 5.75
$ws.Range("Q286").Value = -1.25
$ws.Range("R286").Value = 1.975
$ws.Range("S286").Value = 1.825
$ws.Range("T286").Value = 2.75
$ws.Range("U286").Value = 1.85
$ws.Range("V286").Value = 1.95
$ws.Range("W286").Value = 0.444
$ws.Range("X286").Value = -1
$ws.Range("Y286").Value = -1
$ws.Range("Z286").Value = -0.5
$ws.Range("AA286").Value = 0.4125
$ws.Range("AB286").Value = 0.425
$ws.Range("AC286").Value = -0.5

# Row 294 (new) - copy formatting (bold/border style for id column, date format
# for the date column) from an existing row before overwriting with new values
$ws.Range("A282").Copy($ws.Range("A294"))
$ws.Range("E282").Copy($ws.Range("E294"))
$excel.CutCopyMode = $false
$ws.Range("A294").Value = 292
$ws.Range("B294").Value = 6989634
$ws.Range("C294").Value = 'Serbia Prva Liga'
$ws.Range("D294").Value = 'Serbia Prva Liga'
$ws.Range("E294").Value = 45347.375
$ws.Range("F294").Value = 'FK Graficar Beograd'
$ws.Range("G294").Value = 'FK Indija'
$ws.Range("K294").Value = 2.5
$ws.Range("L294").Value = 3
$ws.Range("M294").Value = 2.625
$ws.Range("N294").Value = 2.15
$ws.Range("O294").Value = 3.1
$ws.Range("P294").Value = 3.1
$ws.Range("Q294").Value = -0.25
$ws.Range("R294").Value = 1.9
$ws.Range("S294").Value = 1.9
$ws.Range("T294").Value = 2.25
$ws.Range("U294").Value = 1.875
$ws.Range("V294").Value = 1.925
$ws.Range("W294").Value = 0
$ws.Range("X294").Value = 0
$ws.Range("Y294").Value = 0
$ws.Range("Z294").Value = 0
$ws.Range("AA294").Value = 0
